$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new cell values for row 2 in the order that matches the target
# shared-string insertion order (D2="Yes" -> 18, E2="err file..." -> 19, C2="dongwang/alphasim_1" -> 20)
$ws.Range("D2").Value = "Yes"
$ws.Range("E2").Value = "err file contains what you expect in the .out file"
$ws.Range("C2").Value = "dongwang/alphasim_1"

# Widen column C
$ws.Range("C:C").ColumnWidth = 22.666666666666668

# Update the active cell selection
$ws.Range("E2").Select()
